# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 969
$wsExhibit.Range("F11").Value = 859
$wsExhibit.Range("F18").Value = 1222
$wsExhibit.Range("F20").Value = 1485
$wsExhibit.Range("F23").Value = 1292
$wsExhibit.Range("F27").Value = 3205

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F5").Value = 59
$wsShow.Range("F11").Value = 24

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 969
$wsAll.Range("F13").Value = 59
$wsAll.Range("F21").Value = 859
$wsAll.Range("F28").Value = 1222
$wsAll.Range("F30").Value = 1485
$wsAll.Range("F33").Value = 1292
$wsAll.Range("F36").Value = 24
$wsAll.Range("F39").Value = 3205
